$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.204.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.545.73"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.996.30"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.072.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.531.27"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.77%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.57"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.46"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0760"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.45%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.42"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.38"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.81%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.64"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.790"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "130.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.25%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.605"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.92"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.29"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.31%  "
